$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 11800
$ws.Range("I70").Value = 16666.666
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 49999.99800000001
$ws.Range("L70").Value = 13500
$ws.Range("M70").Value = -49729.99800000001
$ws.Range("N70").Value = -14040
$ws.Range("H73").Value = 11800
$ws.Range("I73").Value = 16666.666
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 49999.99800000001
$ws.Range("L73").Value = 13500
$ws.Range("M73").Value = -49063.99800000001
$ws.Range("N73").Value = -15372
$ws.Range("H129").Value = 715.72046
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 715.72046
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 2147.16138
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -12147.16138
$ws.Range("H135").Value = 1335.08
$ws.Range("I135").Value = 1523.65
$ws.Range("J135").Value = 580.8
$ws.Range("K135").Value = 13712.85
$ws.Range("L135").Value = 5227.2
$ws.Range("M135").Value = -11177.85
$ws.Range("N135").Value = -10297.2
$ws.Range("H136").Value = 28980
$ws.Range("J136").Value = 28980
$ws.Range("L136").Value = 28980
$ws.Range("N136").Value = -39180
$ws.Range("H137").Value = 10181813
$ws.Range("I137").Value = 15907824
$ws.Range("K137").Value = 47723472
$ws.Range("M137").Value = -47720922
$ws.Range("H138").Value = 2596.46
$ws.Range("I138").Value = 896.9
$ws.Range("J138").Value = 3021.35
$ws.Range("K138").Value = 2690.7
$ws.Range("L138").Value = 9064.05
$ws.Range("M138").Value = 2449.3
$ws.Range("N138").Value = -19344.05

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3735.0908
$ws.Range("I61").Value = 3230.2222
$ws.Range("J61").Value = 6007
$ws.Range("K61").Value = 3230.2222
$ws.Range("L61").Value = 6007
$ws.Range("M61").Value = -3018.2222
$ws.Range("N61").Value = -6431
$ws.Range("H132").Value = 1155.0377
$ws.Range("I132").Value = 712.86365
$ws.Range("J132").Value = 3316.7778
$ws.Range("K132").Value = 2138.59095
$ws.Range("L132").Value = 9950.3334
$ws.Range("M132").Value = 391.4090500000002
$ws.Range("N132").Value = -15010.3334
$ws.Range("H136").Value = 3735.0908
$ws.Range("I136").Value = 3230.2222
$ws.Range("J136").Value = 6007
$ws.Range("K136").Value = 9690.6666
$ws.Range("L136").Value = 18021
$ws.Range("M136").Value = -7140.6666
$ws.Range("N136").Value = -23121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 28747.5
$ws.Range("J141").Value = 34996.668
$ws.Range("L141").Value = 34996.668
$ws.Range("N141").Value = -45356.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3300
$ws.Range("I62").Value = 2600
$ws.Range("K62").Value = 2600
$ws.Range("M62").Value = -1976
$ws.Range("H65").Value = 3300
$ws.Range("I65").Value = 2600
$ws.Range("K65").Value = 13000
$ws.Range("M65").Value = -9880
$ws.Range("H132").Value = 1523.3385
$ws.Range("I132").Value = 1119.0652
$ws.Range("K132").Value = 3357.1956
$ws.Range("M132").Value = -827.1956

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 2650
$ws.Range("I13").Value = 1175
$ws.Range("K13").Value = 3525
$ws.Range("M13").Value = -3357
$ws.Range("H33").Value = 445
$ws.Range("I33").Value = 433.33334
$ws.Range("J33").Value = 480
$ws.Range("K33").Value = 2600.00004
$ws.Range("L33").Value = 2880
$ws.Range("M33").Value = -2317.00004
$ws.Range("N33").Value = -3446
$ws.Range("H44").Value = 425.75
$ws.Range("I44").Value = 401
$ws.Range("J44").Value = 500
$ws.Range("K44").Value = 1203
$ws.Range("L44").Value = 1500
$ws.Range("M44").Value = -805
$ws.Range("N44").Value = -2296
$ws.Range("H57").Value = 4000
$ws.Range("J57").Value = 4000
$ws.Range("L57").Value = 12000
$ws.Range("N57").Value = -13118
$ws.Range("H58").Value = 2749.0908
$ws.Range("J58").Value = 4665
$ws.Range("L58").Value = 13995
$ws.Range("N58").Value = -14251
$ws.Range("H68").Value = 650
$ws.Range("I68").Value = 533.3333
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 1599.9999
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -788.9999
$ws.Range("N68").Value = -4622
$ws.Range("H71").Value = 650
$ws.Range("I71").Value = 533.3333
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 4799.9997
$ws.Range("L71").Value = 9000
$ws.Range("M71").Value = -743.9997000000003
$ws.Range("N71").Value = -17112
$ws.Range("H82").Value = 4041.2
$ws.Range("I82").Value = 3801.5
$ws.Range("J82").Value = 5000
$ws.Range("K82").Value = 11404.5
$ws.Range("L82").Value = 15000
$ws.Range("M82").Value = -10998.5
$ws.Range("N82").Value = -15812
$ws.Range("H85").Value = 4041.2
$ws.Range("I85").Value = 3801.5
$ws.Range("J85").Value = 5000
$ws.Range("K85").Value = 11404.5
$ws.Range("L85").Value = 15000
$ws.Range("M85").Value = -10000.5
$ws.Range("N85").Value = -17808
$ws.Range("H92").Value = 338.8
$ws.Range("I92").Value = 322
$ws.Range("J92").Value = 350
$ws.Range("K92").Value = 966
$ws.Range("L92").Value = 1050
$ws.Range("M92").Value = 282
$ws.Range("N92").Value = -3546
$ws.Range("H97").Value = 749.2857
$ws.Range("I97").Value = 329.875
$ws.Range("J97").Value = 1007.38464
$ws.Range("K97").Value = 989.625
$ws.Range("L97").Value = 3022.15392
$ws.Range("M97").Value = -493.625
$ws.Range("N97").Value = -4014.15392
$ws.Range("H103").Value = 942.6667
$ws.Range("I103").Value = 414
$ws.Range("J103").Value = 2000
$ws.Range("K103").Value = 1242
$ws.Range("L103").Value = 6000
$ws.Range("M103").Value = -363
$ws.Range("N103").Value = -7758
$ws.Range("H109").Value = 3848.5715
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 3848.5715
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 11545.7145
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -13625.7145
$ws.Range("H113").Value = 1642454
$ws.Range("I113").Value = 2873898.2
$ws.Range("J113").Value = 528.1111
$ws.Range("K113").Value = 8621694.600000001
$ws.Range("L113").Value = 1584.3333
$ws.Range("M113").Value = -8619524.600000001
$ws.Range("N113").Value = -5924.3333
$ws.Range("H121").Value = 1207.8108
$ws.Range("I121").Value = 890
$ws.Range("J121").Value = 1257.4688
$ws.Range("K121").Value = 2670
$ws.Range("L121").Value = 3772.4064
$ws.Range("M121").Value = -1360
$ws.Range("N121").Value = -6392.4064
$ws.Range("H122").Value = 703276.3
$ws.Range("I122").Value = 5537.3486
$ws.Range("J122").Value = 2067038.9
$ws.Range("K122").Value = 49836.1374
$ws.Range("L122").Value = 18603350.1
$ws.Range("M122").Value = -47386.1374
$ws.Range("N122").Value = -18608250.1
$ws.Range("H132").Value = 17778008
$ws.Range("I132").Value = 866.6667
$ws.Range("J132").Value = 25396784
$ws.Range("K132").Value = 7800.0003
$ws.Range("L132").Value = 228571056
$ws.Range("M132").Value = -5270.0003
$ws.Range("N132").Value = -228576116

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9323.077
$ws.Range("I70").Value = 9323.077
$ws.Range("K70").Value = 9323.077
$ws.Range("M70").Value = -9053.077
$ws.Range("H73").Value = 9323.077
$ws.Range("I73").Value = 9323.077
$ws.Range("K73").Value = 9323.077
$ws.Range("M73").Value = -8387.077
$ws.Range("H132").Value = 3662.5454
$ws.Range("I132").Value = 3943.3333
$ws.Range("J132").Value = 3060.8572
$ws.Range("K132").Value = 11829.9999
$ws.Range("L132").Value = 9182.5716
$ws.Range("M132").Value = -9299.999899999999
$ws.Range("N132").Value = -14242.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 71430340
$ws.Range("I100").Value = 2071.2
$ws.Range("J100").Value = 250001000
$ws.Range("K100").Value = 2071.2
$ws.Range("L100").Value = 250001000
$ws.Range("M100").Value = -1530.2
$ws.Range("N100").Value = -250002082
$ws.Range("H136").Value = 2658.0322
$ws.Range("I136").Value = 1276.7273
$ws.Range("J136").Value = 3417.75
$ws.Range("K136").Value = 3830.1819
$ws.Range("L136").Value = 10253.25
$ws.Range("M136").Value = -1280.1819
$ws.Range("N136").Value = -15353.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 13347143
$ws.Range("J133").Value = 13347143
$ws.Range("L133").Value = 13347143
$ws.Range("N133").Value = -13357263
$ws.Range("H137").Value = 43770.5
$ws.Range("J137").Value = 43770.5
$ws.Range("L137").Value = 43770.5
$ws.Range("N137").Value = -53970.5
